$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pin Assignment")
$ws.Range("C3").Formula = "=IF(COUNTIF('STK600-RCUC3C0-36 Routing Card'!`$A`$2:`$A`$134,B3),VLOOKUP(B3,'STK600-RCUC3C0-36 Routing Card'!`$A`$2:`$B`$134,2,FALSE),""---"")"
